# Add the new "UNIQUE" column-constraint values to the aspects table.
#
# The source diff shows the "index_letter" (row 3) and "reserve1" (row 6)
# rows gaining a "UNIQUE" marker in column E ("Other 1"), which previously
# held no value. We copy the existing formatting used elsewhere in the
# table (column C's style - Arial 10, no border) onto the target cells so
# the new text matches the surrounding table look, then overwrite the
# value with the actual "UNIQUE" constraint text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Copy($ws.Range("E3"))
$ws.Range("E3").Value = "UNIQUE"

$ws.Range("C6").Copy($ws.Range("E6"))
$ws.Range("E6").Value = "UNIQUE"

# Clear the clipboard marching-ants state left behind by Copy.
$excel.CutCopyMode = 0

# The workbook's last active selection moved to G15.
$ws.Range("G15").Select() | Out-Null

# Best-effort: the author's Excel build/locale renamed the builtin
# "Normal" cell style to "Standard" on save.
try {
    $wb.Styles.Item("Normal").Name = "Standard"
} catch {
}

Write-Output "UNIQUE constraints applied to E3 and E6; selection set to G15."
